$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns with the
# latest scraped values. Values are plain text in the source sheet (e.g.
# "292.48" or "-3.25%"), so we force a Text number format before writing
# them, then restore the Normal style so no stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "292.48"
Set-TextValue $ws.Range("E2") "-3.25%"
Set-TextValue $ws.Range("D3") "31.37"
Set-TextValue $ws.Range("E3") "-1.50%"
Set-TextValue $ws.Range("D4") "4.964"
Set-TextValue $ws.Range("E4") "-1.15%"
Set-TextValue $ws.Range("D5") "0.07379"
Set-TextValue $ws.Range("E5") "-5.40%"
Set-TextValue $ws.Range("D6") "1.809"
Set-TextValue $ws.Range("E6") "-11.54%"
Set-TextValue $ws.Range("D7") "7.660"
Set-TextValue $ws.Range("E7") "-1.55%"
Set-TextValue $ws.Range("D8") "3.750"
Set-TextValue $ws.Range("E8") "-0.84%"
Set-TextValue $ws.Range("D9") "0.9120"
Set-TextValue $ws.Range("E9") "-0.83%"
Set-TextValue $ws.Range("D10") "0.1647"
Set-TextValue $ws.Range("E10") "-5.49%"
Set-TextValue $ws.Range("D11") "0.07628"
Set-TextValue $ws.Range("D12") "0.08190"
Set-TextValue $ws.Range("E12") "-6.94%"
Set-TextValue $ws.Range("D13") "0.02986"
Set-TextValue $ws.Range("E13") "-4.85%"
Set-TextValue $ws.Range("D14") "0.1001"
Set-TextValue $ws.Range("E14") "0.17%"
Set-TextValue $ws.Range("D15") "0.001495"
Set-TextValue $ws.Range("E15") "-0.96%"
Set-TextValue $ws.Range("D16") "0.005645"
Set-TextValue $ws.Range("E16") "-2.93%"
Set-TextValue $ws.Range("D18") "3.469"
Set-TextValue $ws.Range("E18") "0.14%"
Set-TextValue $ws.Range("E19") "-6.16%"
Set-TextValue $ws.Range("D20") "0.3292"
Set-TextValue $ws.Range("E20") "0.00%"
Set-TextValue $ws.Range("D21") "0.1292"
Set-TextValue $ws.Range("E21") "0.00%"
Set-TextValue $ws.Range("D22") "4.333"
Set-TextValue $ws.Range("E22") "3.57%"
Set-TextValue $ws.Range("D23") "0.1976"
Set-TextValue $ws.Range("E23") "9.26%"
Set-TextValue $ws.Range("D24") "0.04483"
Set-TextValue $ws.Range("E24") "-2.74%"
Set-TextValue $ws.Range("D25") "0.001224"
Set-TextValue $ws.Range("E25") "-1.34%"
Set-TextValue $ws.Range("D26") "0.004054"
Set-TextValue $ws.Range("E26") "-9.27%"
Set-TextValue $ws.Range("D27") "0.0001251"
Set-TextValue $ws.Range("E27") "-0.04%"
Set-TextValue $ws.Range("D39") "0.01637"
Set-TextValue $ws.Range("E39") "-5.85%"
Set-TextValue $ws.Range("D40") "0.04401"
Set-TextValue $ws.Range("E40") "-7.45%"
Set-TextValue $ws.Range("D41") "0.007428"
Set-TextValue $ws.Range("E41") "4.35%"
Set-TextValue $ws.Range("E42") "-1.87%"
Set-TextValue $ws.Range("D43") "0.002072"
Set-TextValue $ws.Range("E43") "-1.20%"
Set-TextValue $ws.Range("D44") "0.01112"
Set-TextValue $ws.Range("E44") "2.94%"
Set-TextValue $ws.Range("D45") "0.00006010"
Set-TextValue $ws.Range("E45") "-0.87%"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("E46") "-0.03%"
Set-TextValue $ws.Range("D47") "1.874"
Set-TextValue $ws.Range("E47") "59.44%"
Set-TextValue $ws.Range("D48") "0.003002"
Set-TextValue $ws.Range("E48") "-15.50%"
Set-TextValue $ws.Range("D49") "0.00002101"
Set-TextValue $ws.Range("E49") "-0.03%"
Set-TextValue $ws.Range("D50") "0.0002001"
Set-TextValue $ws.Range("E50") "-0.03%"
